$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-07 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-08 Sunday", 2)

$d.Content.Find.Execute("456×6=2736", $true, $false, $false, $false, $false, $true, 1, $false, "446×5=2230", 2)
$d.Content.Find.Execute("651×9=5859", $true, $false, $false, $false, $false, $true, 1, $false, "238×5=1190", 2)
$d.Content.Find.Execute("540×7=3780", $true, $false, $false, $false, $false, $true, 1, $false, "171×7=1197", 2)
$d.Content.Find.Execute("646×8=5168", $true, $false, $false, $false, $false, $true, 1, $false, "371×6=2226", 2)
$d.Content.Find.Execute("132×4=528", $true, $false, $false, $false, $false, $true, 1, $false, "658×8=5264", 2)

$d.Content.Find.Execute("905×4=3620", $true, $false, $false, $false, $false, $true, 1, $false, "549×6=3294", 2)
$d.Content.Find.Execute("979×3=2937", $true, $false, $false, $false, $false, $true, 1, $false, "784×6=4704", 2)
$d.Content.Find.Execute("197×8=1576", $true, $false, $false, $false, $false, $true, 1, $false, "692×8=5536", 2)
$d.Content.Find.Execute("193×2=386", $true, $false, $false, $false, $false, $true, 1, $false, "728×4=2912", 2)
$d.Content.Find.Execute("225×8=1800", $true, $false, $false, $false, $false, $true, 1, $false, "390×9=3510", 2)

$d.Content.Find.Execute("902×2=1804", $true, $false, $false, $false, $false, $true, 1, $false, "803×3=2409", 2)
$d.Content.Find.Execute("685×8=5480", $true, $false, $false, $false, $false, $true, 1, $false, "696×2=1392", 2)
$d.Content.Find.Execute("231×8=1848", $true, $false, $false, $false, $false, $true, 1, $false, "645×7=4515", 2)
$d.Content.Find.Execute("180×4=720", $true, $false, $false, $false, $false, $true, 1, $false, "672×8=5376", 2)
$d.Content.Find.Execute("739×2=1478", $true, $false, $false, $false, $false, $true, 1, $false, "862×2=1724", 2)

$d.Content.Find.Execute("326×2=652", $true, $false, $false, $false, $false, $true, 1, $false, "386×6=2316", 2)
$d.Content.Find.Execute("255×3=765", $true, $false, $false, $false, $false, $true, 1, $false, "122×6=732", 2)
$d.Content.Find.Execute("274×8=2192", $true, $false, $false, $false, $false, $true, 1, $false, "325×9=2925", 2)
$d.Content.Find.Execute("709×7=4963", $true, $false, $false, $false, $false, $true, 1, $false, "615×7=4305", 2)
$d.Content.Find.Execute("867×6=5202", $true, $false, $false, $false, $false, $true, 1, $false, "124×3=372", 2)

$d.Content.Find.Execute("848×4=3392", $true, $false, $false, $false, $false, $true, 1, $false, "561×6=3366", 2)
$d.Content.Find.Execute("114×7=798", $true, $false, $false, $false, $false, $true, 1, $false, "966×6=5796", 2)
$d.Content.Find.Execute("425×5=2125", $true, $false, $false, $false, $false, $true, 1, $false, "678×2=1356", 2)
$d.Content.Find.Execute("586×4=2344", $true, $false, $false, $false, $false, $true, 1, $false, "394×6=2364", 2)
$d.Content.Find.Execute("864×2=1728", $true, $false, $false, $false, $false, $true, 1, $false, "824×7=5768", 2)
